$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update IMPORTANCIA_CURSOS / MIS_PESOS weights (G2:G6)
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 9

# Update RESULTADOS_X solver assignment column (A2:A876)
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(11, 1).Value = 0
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(23, 1).Value = 0
$ws.Cells.Item(25, 1).Value = 0
$ws.Cells.Item(31, 1).Value = 0
$ws.Cells.Item(45, 1).Value = 1
$ws.Cells.Item(46, 1).Value = 1
$ws.Cells.Item(48, 1).Value = 0
$ws.Cells.Item(51, 1).Value = 0
$ws.Cells.Item(53, 1).Value = 0
$ws.Cells.Item(56, 1).Value = 1
$ws.Cells.Item(60, 1).Value = 1
$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(67, 1).Value = 0
$ws.Cells.Item(69, 1).Value = 0
$ws.Cells.Item(79, 1).Value = 0
$ws.Cells.Item(80, 1).Value = 0
$ws.Cells.Item(83, 1).Value = 1
$ws.Cells.Item(87, 1).Value = 0
$ws.Cells.Item(91, 1).Value = 0
$ws.Cells.Item(94, 1).Value = 0
$ws.Cells.Item(96, 1).Value = 1
$ws.Cells.Item(98, 1).Value = 1
$ws.Cells.Item(99, 1).Value = 0
$ws.Cells.Item(103, 1).Value = 1
$ws.Cells.Item(104, 1).Value = 1
$ws.Cells.Item(105, 1).Value = 1
$ws.Cells.Item(113, 1).Value = 1
$ws.Cells.Item(120, 1).Value = 0
$ws.Cells.Item(140, 1).Value = 0
$ws.Cells.Item(141, 1).Value = 1
$ws.Cells.Item(146, 1).Value = 0
$ws.Cells.Item(148, 1).Value = 0
$ws.Cells.Item(149, 1).Value = 1
$ws.Cells.Item(152, 1).Value = 0
$ws.Cells.Item(155, 1).Value = 1
$ws.Cells.Item(156, 1).Value = 1
$ws.Cells.Item(163, 1).Value = 1
$ws.Cells.Item(164, 1).Value = 1
$ws.Cells.Item(166, 1).Value = 0
$ws.Cells.Item(168, 1).Value = 0
$ws.Cells.Item(171, 1).Value = 1
$ws.Cells.Item(172, 1).Value = 1
$ws.Cells.Item(173, 1).Value = 0
$ws.Cells.Item(176, 1).Value = 0
$ws.Cells.Item(179, 1).Value = 1
$ws.Cells.Item(181, 1).Value = 0
$ws.Cells.Item(184, 1).Value = 0
$ws.Cells.Item(186, 1).Value = 1
$ws.Cells.Item(187, 1).Value = 0
$ws.Cells.Item(191, 1).Value = 0
$ws.Cells.Item(193, 1).Value = 0
$ws.Cells.Item(199, 1).Value = 1
$ws.Cells.Item(201, 1).Value = 1
$ws.Cells.Item(207, 1).Value = 1
$ws.Cells.Item(214, 1).Value = 0
$ws.Cells.Item(219, 1).Value = 1
$ws.Cells.Item(221, 1).Value = 0
$ws.Cells.Item(223, 1).Value = 1
$ws.Cells.Item(226, 1).Value = 1
$ws.Cells.Item(233, 1).Value = 1
$ws.Cells.Item(234, 1).Value = 0
$ws.Cells.Item(240, 1).Value = 0
$ws.Cells.Item(241, 1).Value = 1
$ws.Cells.Item(243, 1).Value = 0
$ws.Cells.Item(247, 1).Value = 0
$ws.Cells.Item(252, 1).Value = 0
$ws.Cells.Item(257, 1).Value = 1
$ws.Cells.Item(264, 1).Value = 1
$ws.Cells.Item(265, 1).Value = 1
$ws.Cells.Item(268, 1).Value = 0
$ws.Cells.Item(270, 1).Value = 1
$ws.Cells.Item(272, 1).Value = 1
$ws.Cells.Item(276, 1).Value = 0
$ws.Cells.Item(277, 1).Value = 0
$ws.Cells.Item(286, 1).Value = 1
$ws.Cells.Item(293, 1).Value = 0
$ws.Cells.Item(295, 1).Value = 1
$ws.Cells.Item(300, 1).Value = 0
$ws.Cells.Item(301, 1).Value = 1
$ws.Cells.Item(305, 1).Value = 0
$ws.Cells.Item(306, 1).Value = 0
$ws.Cells.Item(313, 1).Value = 1
$ws.Cells.Item(314, 1).Value = 1
$ws.Cells.Item(316, 1).Value = 0
$ws.Cells.Item(317, 1).Value = 1
$ws.Cells.Item(322, 1).Value = 1
$ws.Cells.Item(330, 1).Value = 0
$ws.Cells.Item(333, 1).Value = 1
$ws.Cells.Item(334, 1).Value = 0
$ws.Cells.Item(336, 1).Value = 0
$ws.Cells.Item(342, 1).Value = 0
$ws.Cells.Item(345, 1).Value = 1
$ws.Cells.Item(349, 1).Value = 0
$ws.Cells.Item(351, 1).Value = 1
$ws.Cells.Item(354, 1).Value = 0
$ws.Cells.Item(359, 1).Value = 1
$ws.Cells.Item(366, 1).Value = 1
$ws.Cells.Item(367, 1).Value = 0
$ws.Cells.Item(373, 1).Value = 1
$ws.Cells.Item(374, 1).Value = 0
$ws.Cells.Item(375, 1).Value = 1
$ws.Cells.Item(376, 1).Value = 0
$ws.Cells.Item(381, 1).Value = 1
$ws.Cells.Item(382, 1).Value = 0
$ws.Cells.Item(394, 1).Value = 0
$ws.Cells.Item(395, 1).Value = 0
$ws.Cells.Item(397, 1).Value = 0
$ws.Cells.Item(403, 1).Value = 1
$ws.Cells.Item(407, 1).Value = 1
$ws.Cells.Item(408, 1).Value = 0
$ws.Cells.Item(409, 1).Value = 1
$ws.Cells.Item(410, 1).Value = 0
$ws.Cells.Item(415, 1).Value = 1
$ws.Cells.Item(416, 1).Value = 0
$ws.Cells.Item(417, 1).Value = 1
$ws.Cells.Item(419, 1).Value = 1
$ws.Cells.Item(423, 1).Value = 0
$ws.Cells.Item(425, 1).Value = 1
$ws.Cells.Item(426, 1).Value = 1
$ws.Cells.Item(427, 1).Value = 1
$ws.Cells.Item(436, 1).Value = 0
$ws.Cells.Item(441, 1).Value = 1
$ws.Cells.Item(448, 1).Value = 0
$ws.Cells.Item(449, 1).Value = 1
$ws.Cells.Item(450, 1).Value = 0
$ws.Cells.Item(455, 1).Value = 0
$ws.Cells.Item(461, 1).Value = 0
$ws.Cells.Item(463, 1).Value = 0
$ws.Cells.Item(467, 1).Value = 1
$ws.Cells.Item(468, 1).Value = 1
$ws.Cells.Item(473, 1).Value = 0
$ws.Cells.Item(475, 1).Value = 1
$ws.Cells.Item(476, 1).Value = 0
$ws.Cells.Item(481, 1).Value = 1
$ws.Cells.Item(489, 1).Value = 0
$ws.Cells.Item(490, 1).Value = 1
$ws.Cells.Item(493, 1).Value = 1
$ws.Cells.Item(494, 1).Value = 1
$ws.Cells.Item(495, 1).Value = 0
$ws.Cells.Item(497, 1).Value = 0
$ws.Cells.Item(498, 1).Value = 1
$ws.Cells.Item(500, 1).Value = 1
$ws.Cells.Item(503, 1).Value = 0
$ws.Cells.Item(507, 1).Value = 1
$ws.Cells.Item(510, 1).Value = 0
$ws.Cells.Item(512, 1).Value = 0
$ws.Cells.Item(518, 1).Value = 1
$ws.Cells.Item(519, 1).Value = 0
$ws.Cells.Item(527, 1).Value = 0
$ws.Cells.Item(528, 1).Value = 0
$ws.Cells.Item(529, 1).Value = 0
$ws.Cells.Item(530, 1).Value = 0
$ws.Cells.Item(531, 1).Value = 0
$ws.Cells.Item(543, 1).Value = 1
$ws.Cells.Item(545, 1).Value = 0
$ws.Cells.Item(551, 1).Value = 1
$ws.Cells.Item(552, 1).Value = 1
$ws.Cells.Item(556, 1).Value = 1
$ws.Cells.Item(558, 1).Value = 1
$ws.Cells.Item(559, 1).Value = 1
$ws.Cells.Item(562, 1).Value = 1
$ws.Cells.Item(563, 1).Value = 1
$ws.Cells.Item(564, 1).Value = 1
$ws.Cells.Item(565, 1).Value = 1
$ws.Cells.Item(569, 1).Value = 0
$ws.Cells.Item(570, 1).Value = 0
$ws.Cells.Item(571, 1).Value = 0
$ws.Cells.Item(573, 1).Value = 0
$ws.Cells.Item(575, 1).Value = 0
$ws.Cells.Item(577, 1).Value = 1
$ws.Cells.Item(580, 1).Value = 1
$ws.Cells.Item(582, 1).Value = 1
$ws.Cells.Item(590, 1).Value = 0
$ws.Cells.Item(593, 1).Value = 0
$ws.Cells.Item(602, 1).Value = 0
$ws.Cells.Item(603, 1).Value = 0
$ws.Cells.Item(604, 1).Value = 1
$ws.Cells.Item(605, 1).Value = 1
$ws.Cells.Item(607, 1).Value = 1
$ws.Cells.Item(612, 1).Value = 0
$ws.Cells.Item(613, 1).Value = 0
$ws.Cells.Item(614, 1).Value = 1
$ws.Cells.Item(623, 1).Value = 1
$ws.Cells.Item(624, 1).Value = 1
$ws.Cells.Item(626, 1).Value = 0
$ws.Cells.Item(630, 1).Value = 0
$ws.Cells.Item(638, 1).Value = 1
$ws.Cells.Item(641, 1).Value = 1
$ws.Cells.Item(643, 1).Value = 1
$ws.Cells.Item(644, 1).Value = 1
$ws.Cells.Item(645, 1).Value = 1
$ws.Cells.Item(649, 1).Value = 0
$ws.Cells.Item(652, 1).Value = 0
$ws.Cells.Item(653, 1).Value = 0
$ws.Cells.Item(654, 1).Value = 0
$ws.Cells.Item(656, 1).Value = 0
$ws.Cells.Item(658, 1).Value = 0
$ws.Cells.Item(659, 1).Value = 0
$ws.Cells.Item(662, 1).Value = 1
$ws.Cells.Item(665, 1).Value = 1
$ws.Cells.Item(671, 1).Value = 1
$ws.Cells.Item(672, 1).Value = 1
$ws.Cells.Item(677, 1).Value = 0
$ws.Cells.Item(679, 1).Value = 0
$ws.Cells.Item(688, 1).Value = 1
$ws.Cells.Item(689, 1).Value = 1
$ws.Cells.Item(692, 1).Value = 0
$ws.Cells.Item(695, 1).Value = 1
$ws.Cells.Item(697, 1).Value = 0
$ws.Cells.Item(699, 1).Value = 0
$ws.Cells.Item(709, 1).Value = 0
$ws.Cells.Item(710, 1).Value = 0
$ws.Cells.Item(711, 1).Value = 0
$ws.Cells.Item(712, 1).Value = 0
$ws.Cells.Item(718, 1).Value = 0
$ws.Cells.Item(723, 1).Value = 1
$ws.Cells.Item(724, 1).Value = 1
$ws.Cells.Item(725, 1).Value = 1
$ws.Cells.Item(730, 1).Value = 1
$ws.Cells.Item(732, 1).Value = 1
$ws.Cells.Item(743, 1).Value = 1
$ws.Cells.Item(744, 1).Value = 1
$ws.Cells.Item(745, 1).Value = 1
$ws.Cells.Item(746, 1).Value = 1
$ws.Cells.Item(749, 1).Value = 0
$ws.Cells.Item(751, 1).Value = 1
$ws.Cells.Item(753, 1).Value = 1
$ws.Cells.Item(756, 1).Value = 1
$ws.Cells.Item(758, 1).Value = 0
$ws.Cells.Item(759, 1).Value = 0
$ws.Cells.Item(760, 1).Value = 0
$ws.Cells.Item(766, 1).Value = 0
$ws.Cells.Item(767, 1).Value = 0
$ws.Cells.Item(769, 1).Value = 0
$ws.Cells.Item(772, 1).Value = 0
$ws.Cells.Item(773, 1).Value = 0
$ws.Cells.Item(774, 1).Value = 0
$ws.Cells.Item(783, 1).Value = 1
$ws.Cells.Item(784, 1).Value = 1
$ws.Cells.Item(785, 1).Value = 1
$ws.Cells.Item(786, 1).Value = 0
$ws.Cells.Item(790, 1).Value = 1
$ws.Cells.Item(796, 1).Value = 1
$ws.Cells.Item(797, 1).Value = 1
$ws.Cells.Item(800, 1).Value = 0
$ws.Cells.Item(806, 1).Value = 0
$ws.Cells.Item(813, 1).Value = 0
$ws.Cells.Item(817, 1).Value = 1
$ws.Cells.Item(825, 1).Value = 0
$ws.Cells.Item(826, 1).Value = 0
$ws.Cells.Item(832, 1).Value = 0
$ws.Cells.Item(838, 1).Value = 1
$ws.Cells.Item(839, 1).Value = 1
$ws.Cells.Item(841, 1).Value = 1
$ws.Cells.Item(842, 1).Value = 1
$ws.Cells.Item(843, 1).Value = 1
$ws.Cells.Item(844, 1).Value = 1
$ws.Cells.Item(853, 1).Value = 0
$ws.Cells.Item(855, 1).Value = 0
$ws.Cells.Item(866, 1).Value = 0
$ws.Cells.Item(871, 1).Value = 1
$ws.Cells.Item(873, 1).Value = 0

# Update active selection
$ws.Range("H7").Select()
